# Week 16 logged + season sim from Week 17:
# a new player ("E.Wolf") joined the tracked roster, so a new column is
# inserted right before the "A.Trautman" column (column W) on both the
# "Rushing" and "Receiving" sheets, pushing A.Trautman/J.Johnson/
# G.Griffin/N.Vannett one column to the right. The new column gets the
# same header style as the rest of row 1, and an "n" placeholder in the
# data row, matching its neighbours.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Insert a new blank column at W (before A.Trautman), shifting
    # W:Z -> X:AA and inheriting the existing header/data formatting.
    $ws.Columns("W:W").Insert()

    # New player header (row 1) and placeholder data (row 2).
    $ws.Range("W1").Value = "E.Wolf"
    $ws.Range("W2").Value = "n"
}
